# fix templates for dictionary uploads
#
# On the "Variable values" sheet, the header in column C was "name" and
# should instead read "variable" (matching the new molgenis-emx2 upload
# template). Update the cell, widen the column so the new header still
# fits, and leave the cell selection on C2 - exactly like a user would
# after typing the new header and pressing Enter.

$wb = $excel.ActiveWorkbook

$dsSheet = $wb.Worksheets.Item(1)        # "Datasets" - stays the active tab
$valuesSheet = $wb.Worksheets.Item(3)    # "Variable values"

# Rename header cell C1 from "name" to "variable"
$valuesSheet.Cells.Item(1, 3).Value = "variable"

# Widen column C a bit so the longer header text keeps fitting
$valuesSheet.Columns.Item(3).ColumnWidth = 6.03

# Leave the selection on C2 (just below the edited header), then restore
# the originally active sheet/tab
[void]$valuesSheet.Range("C2").Select()
[void]$dsSheet.Select()
